$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.201.42"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "'3.578.27"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'192.62"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").Value = "'573.45"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("D7").Value = "'3.571.37"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").Value = "'0.616"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'0.677"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").Value = "'55.64"
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "'9.88"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "'4.150.11"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("D16").Value = "'3.573.37"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "'67.062.82"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'12.24"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "'18.26"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "'403.60"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").Value = "'4.18"
$ws.Range("E23").Value = "  -4.35%  "
$ws.Range("D24").Value = "'12.16"
$ws.Range("E24").Value = "  +8.83%  "
$ws.Range("D25").Value = "'85.80"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").Value = "'2.92"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'12.56"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "'6.10"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").Value = "'3.71"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").Value = "'7.85"
$ws.Range("E30").Value = "  +6.78%  "
$ws.Range("D31").Value = "'8.99"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").Value = "'31.24"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "'644.17"
$ws.Range("E33").Value = "  +6.20%  "
$ws.Range("D34").Value = "'12.15"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").Value = "'63.79"
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("D37").Value = "'42.40"
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("D38").Value = "'0.409"
$ws.Range("E38").Value = "  +4.53%  "
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'0.0₃0772"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'3.173.88"
$ws.Range("E41").Value = "  +13.53%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "'3.10"
$ws.Range("E42").Value = "  +7.80%  "
$ws.Range("D43").Value = "'0.134"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.73"
$ws.Range("E44").Value = "  +8.64%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "'0.0416"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.130"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.10"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'143.24"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'8.58"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").Value = "'2.53"
$ws.Range("E51").Value = "  -3.42%  "
